# BOT; UPDATE DATA
# Adds one more day of data (2020-05-01) to the "相談件数" (consultation
# counts) table, pushing the trailing footnote row down by one row, and
# updates the sheet's dimension/selection to reflect the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 97. This shifts the old row 97 (the footnote
# "※4/8より健康相談窓口と帰国者・接触者相談センターを統合") down to row 98,
# and the new row 97 inherits the number formatting from row 96 above it
# (date format in A, "0_);[Red](0)" in B/C, general/right-aligned in D/E).
$ws.Rows.Item(97).Insert()

# Populate the new row with the 2020-05-01 figures.
$ws.Cells.Item(97, 1).Value = [datetime]"2020-05-01"   # date (serial 43952)
$ws.Cells.Item(97, 2).Value = 407                        # 保健所・保健センター（日別）
$ws.Cells.Item(97, 3).Value = 32436                       # 保健所・保健センター（累計）
$ws.Cells.Item(97, 4).Value = 109                         # 専用健康相談窓口（日別）
$ws.Cells.Item(97, 5).Value = 6958                        # 専用健康相談窓口（累計）

# The data range (and the sheet's used range) now ends at row 98; move the
# selection to the new bottom-right corner of the table to match.
$ws.Range("E98").Select()
